$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in phone numbers for the two data rows
$ws.Range("F2").Value = "0962772733"
$ws.Range("F3").Value = "0962772733"

# Fill in employee name "Trung" in A3
$ws.Range("A3").Value = "Trung"

# Update the active selection to A3
$ws.Range("A3").Select()
